# ---------------------------------------------------------------------------
# Applies the "Add files via upload" edit described by the diff:
#   1. Adds a new centered/bold "19.06.2023" paragraph right after "DAY 01".
#   2. Drops the (now resolved) spell-check proofErr bookmarks around the
#      seven CNN-architecture names (LeNet, AlexNet, VGGNet, GoogLeNet,
#      MobileNet, ResNet, DenseNet).
#   3. Drops the stale <w:lastRenderedPageBreak/> marker in front of the
#      "Project: Traffic Sign Detection" heading.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$wdCollapseEnd = 0

function Find-ParagraphByExactText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs($i)
        $t = $para.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $para
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "DAY 01" -> insert a new paragraph "19.06.2023" right after it, copying
#    the same paragraph/run formatting (centered, bold, sz 28).
# ---------------------------------------------------------------------------
$dayPara = Find-ParagraphByExactText("DAY 01")
if ($dayPara -ne $null) {
    $dayPara.Range.InsertParagraphAfter()
    $datePara = $dayPara.Next()
    $datePara.Range.Text = "19.06.2023"
}

# ---------------------------------------------------------------------------
# 2) Remove the <w:proofErr .../> spellStart/spellEnd pair wrapping each of
#    these single-word runs - rebuild the paragraph via InsertXML without
#    the proofErr markers, keeping every other property identical.
# ---------------------------------------------------------------------------
function Remove-ProofErrParagraph($text, $spacingAfter0) {
    $para = Find-ParagraphByExactText($text)
    if ($para -eq $null) {
        return
    }
    $spacing = ""
    if ($spacingAfter0) {
        $spacing = '<w:spacing w:after="0"/>'
    }
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`n" + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + $spacing + '<w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>' + $text + '</w:t></w:r>' + `
        '</w:p>' + `
        '</w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($xml)
}

Remove-ProofErrParagraph "LeNet" $true
Remove-ProofErrParagraph "AlexNet" $true
Remove-ProofErrParagraph "VGGNet" $true
Remove-ProofErrParagraph "GoogLeNet" $true
Remove-ProofErrParagraph "MobileNet" $true
Remove-ProofErrParagraph "ResNet" $true
Remove-ProofErrParagraph "DenseNet" $false

# ---------------------------------------------------------------------------
# 3) Drop the stray <w:lastRenderedPageBreak/> before "Project: Traffic Sign
#    Detection" (same rebuild-via-InsertXML technique).
# ---------------------------------------------------------------------------
$projPara = Find-ParagraphByExactText("Project: Traffic Sign Detection")
if ($projPara -ne $null) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`n" + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
        '<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Project: Traffic Sign Detection</w:t></w:r>' + `
        '</w:p>' + `
        '</w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
    $projPara.Range.InsertXML($xml)
}

Write-Output "edit complete"
